# Updated cryptos list on Mon Nov 11 07:58:57 UTC 2024 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row, and
# re-ranks three coins (Stacks/WhiteBITCoin/dogwifhat at rows 42-44, and
# InjectiveProtocol/Mantle at rows 50-51) to match their new positions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D are plain text (e.g. "208.39", "81.155.89") in the
# source sheet. Excel auto-coerces a bare numeric-looking string typed into
# Value to a number, so force the cell to Text first, then restore the
# default "Normal" style afterwards to avoid leaving stray number formatting
# behind (multi-dot values like "81.155.89" are never auto-numeric, but the
# same NumberFormat/Style bracket is applied uniformly for all D-column
# writes for simplicity).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '81.155.89'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.94%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.142.51'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.38%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '208.39'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '616.53'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.282'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +24.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -0.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.138.89'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.40%  '
$ws.Range('E11').Value = '  -2.04%  '
$ws.Range('E12').Value = '  +11.75%  '
$ws.Range('E13').Value = '  -0.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.27'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.718.47'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.49%  '
$ws.Range('E16').Value = '  -0.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '80.888.97'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.70%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.145.69'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.17'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +9.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.82'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '430.06'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.95'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.08'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.17'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.73%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.19'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +9.33%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.300.90'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '75.65'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.19%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.83'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.73%  '
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('E30').Value = '  +5.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.998'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.25%  '
$ws.Range('E32').Value = '  +0.89%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '566.80'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +9.94%  '
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('E35').Value = '  +14.77%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.137'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +8.74%  '
$ws.Range('E37').Value = '  -0.82%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '22.53'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.86%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.998'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.11'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +12.95%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.405'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.75%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.04'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +15.03%  '
$ws.Range('B43').Value = 'WhiteBITCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '20.73'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.73%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.03'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +21.85%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '158.99'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.98%  '
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '186.58'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.98%  '
$ws.Range('E48').Value = '  +4.60%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.31'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.38%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '25.98'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.47%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.762'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.17%  '
